$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.801.02'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '3.052.19'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'556.50"
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').Value = "'141.87"
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.051.65'
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('E11').Value = '  -13.18%  '
$ws.Range('E12').Value = '  +1.86%  '
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = "'35.13"
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '3.548.65'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').Value = '63.818.37'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Value = '3.050.11'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').Value = "'488.26"
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').Value = "'14.15"
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').Value = "'0.682"
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = "'14.42"
$ws.Range('E23').Value = '  +6.60%  '
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').Value = "'82.52"
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = "'8.10"
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('D29').Value = "'2.03"
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = "'26.27"
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').Value = "'5.67"
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('D35').Value = "'6.19"
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('D36').Value = "'55.25"
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = "'0.0408"
$ws.Range('E37').Value = '  -0.86%  '
$ws.Range('D38').Value = "'441.30"
$ws.Range('E38').Value = '  -6.73%  '
$ws.Range('D39').Value = "'0.0814"
$ws.Range('E39').Value = '  -2.30%  '
$ws.Range('D40').Value = '3.015.75'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = "'2.76"
$ws.Range('E41').Value = '  -6.39%  '
$ws.Range('D42').Value = "'8.31"
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('E44').Value = '  +4.26%  '
$ws.Range('D45').Value = "'27.67"
$ws.Range('E45').Value = '  -2.42%  '
$ws.Range('E46').Value = '  +3.37%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').Value = "'118.01"
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').Value = '0.0₃0512'
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('D51').Value = "'2.09"
$ws.Range('E51').Value = '  -0.24%  '
